$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing ingredient names (case / wording changes), in the same
# order the new shared-string entries were introduced by the author.
$ws.Range("A2").Value = "flour"

$ws.Range("A4").Value = "oil"

# Add a new ingredient row
$ws.Range("A14").Value = "vanilla"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "bottle"

$ws.Range("A10").Value = "white sugar"

$ws.Range("A3").Value = "salt"

# Restore selection to A3 as in the saved file
$ws.Range("A3").Select()
